# Updated symbol list on Sat Dec 31 09:21:48 UTC 2022 with GitHub Actions
# Applies the latest coinranking.com scrape results to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price-only updates (column D) ---------------------------------------
# Force text formatting so values like "3.551" / "0.00009600" keep their
# exact textual representation (trailing zeros, etc.) instead of being
# reinterpreted as floating point numbers.
$priceCells = @("D2","D3","D4","D5","D6","D8","D9","D15","D16","D18","D20","D22","D23","D24","D25","D26","D27","D28","D40","D44","D45","D47","D48","D49","D50")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "245.76"
$ws.Range("D3").Value = "25.91"
$ws.Range("D4").Value = "5.115"
$ws.Range("D5").Value = "0.05589"
$ws.Range("D6").Value = "6.483"
$ws.Range("D8").Value = "0.8164"
$ws.Range("D9").Value = "0.8484"
$ws.Range("D15").Value = "0.006123"
$ws.Range("D16").Value = "3.551"
$ws.Range("D18").Value = "0.3134"
$ws.Range("D20").Value = "0.03219"
$ws.Range("D22").Value = "3.739"
$ws.Range("D23").Value = "0.04704"
$ws.Range("D24").Value = "0.1375"
$ws.Range("D25").Value = "0.001250"
$ws.Range("D26").Value = "0.004606"
$ws.Range("D27").Value = "0.00009600"
$ws.Range("D28").Value = "0.0001390"
$ws.Range("D40").Value = "0.03655"
$ws.Range("D44").Value = "0.007883"
$ws.Range("D45").Value = "0.00005314"
$ws.Range("D47").Value = "0.1335"
$ws.Range("D48").Value = "0.002048"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("D50").Value = "0.0002000"

# --- Row 10: One -> WazirX -------------------------------------------------
$ws.Range("D10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1337"
$ws.Range("E10").Value = "9WazirXWRX"

# --- Row 11: WazirX -> BitrueCoin ------------------------------------------
$ws.Range("D11").NumberFormat = "@"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "0.02853"
$ws.Range("E11").Value = "10BitrueCoinBTR"

# --- Row 12: BitrueCoin -> BitMartToken ------------------------------------
$ws.Range("D12").NumberFormat = "@"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "0.09407"
$ws.Range("E12").Value = "11BitMartTokenBMX"

# --- Row 13: BitMartToken -> BitForexToken ---------------------------------
$ws.Range("D13").NumberFormat = "@"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "0.001518"
$ws.Range("E13").Value = "12BitForexTokenBF"

# --- Row 14: BitForexToken -> One ------------------------------------------
$ws.Range("D14").NumberFormat = "@"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "0.0006001"
$ws.Range("E14").Value = "13OneONE"

# --- Row 41: KickToken -> BKEXToken ----------------------------------------
$ws.Range("D41").NumberFormat = "@"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "0.1351"
$ws.Range("E41").Value = "40BKEXTokenBKKBestin24h"

# --- Row 42: BKEXToken -> CEJI ----------------------------------------------
$ws.Range("D42").NumberFormat = "@"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.002483"
$ws.Range("E42").Value = "41CEJICEJI"

# --- Row 43: CEJI -> KickToken ----------------------------------------------
$ws.Range("D43").NumberFormat = "@"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "0.003389"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

# --- Row 47: Worst-in-24h badge removed -------------------------------------
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
